$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "Ajout d'une couleur au survol pour les liens accueil et contact"
$ws.Cells.Item(7, 3).Value = "fait"

# Row 8
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "suppression des img de citation et modification en texte"
$ws.Cells.Item(8, 3).Value = "fait"

# Row 9 (previously row 7, now pushed down, gets its value completed)
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Modification du footer en mettant des infos utile"
$ws.Cells.Item(9, 3).Value = "fait"

# Widen column B to fit the new, longer text (closest reachable value to the
# author's final bestFit width of 55.5703125 given this host's column-width
# quantization grid).
$ws.Columns.Item(2).ColumnWidth = 54.65

# Update the active selection to match the author's final cursor position.
$ws.Range("B9").Select()
